# Auto-generated edit script: update crypto price/volume table (columns D & E)
# for rows 2-51, plus a row-49/50 content swap (B/C/D/E), matching the
# "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to keep a literal text value (e.g. "318.03" or
    # "45.392.63") instead of Excel auto-coercing it into a number, while
    # restoring the cell style afterwards so no stray number format sticks.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "45.392.63"
Set-TextValue $ws.Range("E2") "  +5.85%  "
Set-TextValue $ws.Range("D3") "2.380.26"
Set-TextValue $ws.Range("E3") "  +3.70%  "
Set-TextValue $ws.Range("E4") "  +0.36%  "
Set-TextValue $ws.Range("D5") "111.97"
Set-TextValue $ws.Range("E5") "  +7.21%  "
Set-TextValue $ws.Range("D6") "318.03"
Set-TextValue $ws.Range("E6") "  +2.02%  "
Set-TextValue $ws.Range("E7") "  +2.34%  "
Set-TextValue $ws.Range("E8") "  -0.18%  "
Set-TextValue $ws.Range("E9") "  +5.12%  "
Set-TextValue $ws.Range("D10") "42.34"
Set-TextValue $ws.Range("E10") "  +8.07%  "
Set-TextValue $ws.Range("D11") "0.0930"
Set-TextValue $ws.Range("E11") "  +3.03%  "
Set-TextValue $ws.Range("E12") "  +5.28%  "
Set-TextValue $ws.Range("E13") "  +3.08%  "
Set-TextValue $ws.Range("D14") "0.108"
Set-TextValue $ws.Range("E14") "  +0.25%  "
Set-TextValue $ws.Range("D15") "15.83"
Set-TextValue $ws.Range("E15") "  +4.63%  "
Set-TextValue $ws.Range("D16") "2.743.29"
Set-TextValue $ws.Range("E16") "  +3.73%  "
Set-TextValue $ws.Range("D17") "2.389.42"
Set-TextValue $ws.Range("E17") "  +3.25%  "
Set-TextValue $ws.Range("D18") "45.365.31"
Set-TextValue $ws.Range("E18") "  +6.72%  "
Set-TextValue $ws.Range("D19") "7.63"
Set-TextValue $ws.Range("E19") "  +4.51%  "
Set-TextValue $ws.Range("E20") "  +3.51%  "
Set-TextValue $ws.Range("D21") "13.11"
Set-TextValue $ws.Range("E21") "  -2.64%  "
Set-TextValue $ws.Range("D22") "75.26"
Set-TextValue $ws.Range("E22") "  +2.72%  "
Set-TextValue $ws.Range("E23") "  +2.91%  "
Set-TextValue $ws.Range("D24") "269.21"
Set-TextValue $ws.Range("E24") "  +2.20%  "
Set-TextValue $ws.Range("D25") "2.34"
Set-TextValue $ws.Range("E25") "  +7.28%  "
Set-TextValue $ws.Range("E26") "  -0.38%  "
Set-TextValue $ws.Range("D27") "7.71"
Set-TextValue $ws.Range("E27") "  +9.01%  "
Set-TextValue $ws.Range("D28") "11.31"
Set-TextValue $ws.Range("E28") "  +5.44%  "
Set-TextValue $ws.Range("E29") "  +0.10%  "
Set-TextValue $ws.Range("D30") "39.31"
Set-TextValue $ws.Range("E30") "  +9.88%  "
Set-TextValue $ws.Range("D31") "22.90"
Set-TextValue $ws.Range("E31") "  +2.45%  "
Set-TextValue $ws.Range("D32") "0.0935"
Set-TextValue $ws.Range("E32") "  +8.37%  "
Set-TextValue $ws.Range("D33") "169.69"
Set-TextValue $ws.Range("E33") "  +2.87%  "
Set-TextValue $ws.Range("E34") "  +15.82%  "
Set-TextValue $ws.Range("E35") "  +2.17%  "
Set-TextValue $ws.Range("E36") "  +3.71%  "
Set-TextValue $ws.Range("E37") "  +7.50%  "
Set-TextValue $ws.Range("D38") "3.06"
Set-TextValue $ws.Range("E38") "  +12.04%  "
Set-TextValue $ws.Range("D39") "0.0369"
Set-TextValue $ws.Range("E39") "  +5.27%  "
Set-TextValue $ws.Range("E40") "  +5.35%  "
Set-TextValue $ws.Range("D41") "1.75"
Set-TextValue $ws.Range("E41") "  +9.78%  "
Set-TextValue $ws.Range("D42") "105.86"
Set-TextValue $ws.Range("E42") "  +6.47%  "
Set-TextValue $ws.Range("E43") "  +15.48%  "
Set-TextValue $ws.Range("D44") "0.241"
Set-TextValue $ws.Range("E44") "  +6.52%  "
Set-TextValue $ws.Range("D45") "72.16"
Set-TextValue $ws.Range("E45") "  +4.22%  "
Set-TextValue $ws.Range("E46") "  +0.30%  "
Set-TextValue $ws.Range("D47") "119.12"
Set-TextValue $ws.Range("E47") "  +7.39%  "
Set-TextValue $ws.Range("D48") "5.72"
Set-TextValue $ws.Range("E48") "  +10.10%  "
Set-TextValue $ws.Range("B49") "ordi"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextValue $ws.Range("D49") "80.05"
Set-TextValue $ws.Range("E49") "  +0.74%  "
Set-TextValue $ws.Range("B50") "MinaProtocolToken"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
Set-TextValue $ws.Range("D50") "1.63"
Set-TextValue $ws.Range("E50") "  +18.63%  "
Set-TextValue $ws.Range("E51") "  +17.03%  "
